# Add 2022-Q3 data:
#  1. Insert a new worksheet named "2022-Q3" before the existing "2022-Q2" sheet.
#  2. Populate it with the fund holdings table for 2022-Q3.
#  3. Insert a new row at the top of the "总计" (summary) sheet's data table
#     for the 2022-Q3 totals, shifting the existing rows down.
#
# NOTE: worksheet object variables in this COM shim resolve positionally
# (i.e. a captured reference tracks "the sheet currently at index N", not a
# stable identity). After inserting a new sheet, any previously-captured
# reference to a sheet that sat at/after the insertion point must be
# re-fetched by name before being used again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right before "2022-Q2"
# ---------------------------------------------------------------------------
$q2SheetOriginal = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2SheetOriginal)
$q3Sheet.Name = "2022-Q3"

# Re-fetch "2022-Q2" by name now that the sheet collection has shifted --
# $q2SheetOriginal is no longer the "2022-Q2" sheet after the insert.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Header row (row 1): B..H, bold/centered/bordered style matching the other
# quarter sheets (copy formatting from the corresponding header cells of the
# "2022-Q2" sheet).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3Sheet.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
}
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats

# Data rows (row 2..39): A + H are numeric, B..G are text.
$rows = @(
    @(0,  "398051", "中海环保新能源混合",             "22.19", "67.74", "5.93", "1.3159", 4),
    @(1,  "398021", "中海能源策略混合",                 "21.73", "88.03", "3.72", "0.8084", 9),
    @(2,  "200015", "长城优化升级混合A",                "17.65", "88.77", "3.24", "0.5719", 8),
    @(3,  "519087", "新华优选分红混合",                 "10.82", "89.55", "3.93", "0.4252", 9),
    @(4,  "001279", "中海积极增利灵活配置混合",         "6.90",  "83.64", "3.73", "0.2574", 9),
    @(5,  "519156", "新华行业轮换灵活配置混合A",        "5.70",  "94.21", "4.43", "0.2525", 8),
    @(6,  "013274", "长城优化升级混合C",                "6.73",  "88.77", "3.24", "0.2181", 8),
    @(7,  "161039", "富国中证1000指数增强（LOF）A",     "25.41", "84.72", "0.70", "0.1779", 4),
    @(8,  "000800", "华商未来主题混合",                 "4.21",  "72.31", "4.04", "0.1701", 5),
    @(9,  "519975", "长信量化中小盘股票",               "8.16",  "93.71", "1.87", "0.1526", 10),
    @(10, "013916", "中融成长先锋一年持有混合A",        "2.15",  "88.19", "3.62", "0.0778", 8),
    @(11, "013331", "富国中证1000指数增强（LOF）C",     "8.53",  "84.72", "0.70", "0.0597", 4),
    @(12, "014202", "天弘中证1000指数增强C",            "3.69",  "94.06", "1.58", "0.0583", 6),
    @(13, "014201", "天弘中证1000指数增强A",            "3.68",  "94.06", "1.58", "0.0581", 6),
    @(14, "014329", "中融优势产业混合A",                "1.48",  "69.84", "3.85", "0.0570", 7),
    @(15, "005632", "鹏华量化先锋混合",                 "2.71",  "92.57", "2.10", "0.0569", 2),
    @(16, "168207", "中融创业板两年定期开放混合",       "1.02",  "83.05", "5.33", "0.0544", 5),
    @(17, "014571", "东吴安享量化灵活配置混合C",        "0.47",  "90.86", "9.33", "0.0439", 7),
    @(18, "580007", "东吴安享量化灵活配置混合A",        "0.47",  "90.86", "9.33", "0.0439", 7),
    @(19, "011457", "新华行业龙头主题股票",             "0.97",  "94.28", "4.10", "0.0398", 8),
    @(20, "002210", "创金合信量化多因子股票A",          "2.39",  "91.71", "1.28", "0.0306", 5),
    @(21, "013466", "博时智选量化多因子股票C",          "2.28",  "92.38", "1.02", "0.0233", 10),
    @(22, "003865", "创金合信量化多因子股票C",          "0.75",  "91.71", "1.28", "0.0096", 5),
    @(23, "014014", "招商臻选平衡混合A",                "0.33",  "66.43", "2.46", "0.0081", 8),
    @(24, "014015", "招商臻选平衡混合C",                "0.26",  "66.43", "2.46", "0.0064", 8),
    @(25, "320016", "诺安多策略混合",                   "0.17",  "76.95", "3.77", "0.0064", 9),
    @(26, "011731", "国投瑞银安睿混合A",                "1.21",  "26.60", "0.42", "0.0051", 4),
    @(27, "013465", "博时智选量化多因子股票A",          "0.49",  "92.38", "1.02", "0.0050", 10),
    @(28, "015466", "太平中证1000指数增强A",            "0.37",  "92.23", "1.07", "0.0040", 4),
    @(29, "009514", "创金合信同顺创业板精选股票C",      "0.16",  "92.10", "2.41", "0.0039", 7),
    @(30, "014330", "中融优势产业混合C",                "0.10",  "69.84", "3.85", "0.0038", 7),
    @(31, "011732", "国投瑞银安睿混合C",                "0.83",  "26.60", "0.42", "0.0035", 4),
    @(32, "004360", "创金合信量化核心混合C",            "0.21",  "91.03", "1.40", "0.0029", 8),
    @(33, "519157", "新华行业轮换灵活配置混合C",        "0.06",  "94.21", "4.43", "0.0027", 8),
    @(34, "013917", "中融成长先锋一年持有混合C",        "0.07",  "88.19", "3.62", "0.0025", 8),
    @(35, "004359", "创金合信量化核心混合A",            "0.16",  "91.03", "1.40", "0.0022", 8),
    @(36, "009513", "创金合信同顺创业板精选股票A",      "0.09",  "92.10", "2.41", "0.0022", 7),
    @(37, "015467", "太平中证1000指数增强C",            "0.02",  "92.23", "1.07", "0.0002", 4)
)

# Force B:G to be stored as text (not auto-converted to numbers) by setting
# the number format to "@" before assigning the values.
$q3Sheet.Range("B2:G39").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q3Sheet.Range("A$r").Value = $row[0]
    $q3Sheet.Range("B$r").Value = $row[1]
    $q3Sheet.Range("C$r").Value = $row[2]
    $q3Sheet.Range("D$r").Value = $row[3]
    $q3Sheet.Range("E$r").Value = $row[4]
    $q3Sheet.Range("F$r").Value = $row[5]
    $q3Sheet.Range("G$r").Value = $row[6]
    $q3Sheet.Range("H$r").Value = $row[7]
}

# Restore "General" number format on the text columns now that the values
# are already stored as text -- matches the un-styled text cells used by
# the other quarter sheets.
$q3Sheet.Range("B2:G39").Style = "Normal"

# Column A uses the bold/centered/bordered "row index" style, same as the
# other quarter sheets -- copy it from the equivalent cell on "2022-Q2".
$q2Sheet.Range("A2").Copy()
$q3Sheet.Range("A2:A39").PasteSpecial(-4122) # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Insert the new 2022-Q3 total row into the "总计" summary sheet
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-apply the row-index style ("s=2": bold/centered/bordered) to A2, which
# Insert() leaves unstyled, and clear the format stamped onto B2:D2 by the
# insert operation so they match the plain (unstyled) data cells used
# elsewhere in the table.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 38
$totalSheet.Range("D2").Value = 5.02

# Column A is a plain 0-based row index (0,1,2,...) -- Insert() shifts the
# old cell contents down with the row, so A3..A8 currently still hold the
# pre-insert values (0,1,2,3,4,5 instead of 1,2,3,4,5,6). Renumber the whole
# index column explicitly so it reads 0..6 top to bottom.
for ($r = 2; $r -le 8; $r++) {
    $totalSheet.Range("A$r").Value = $r - 2
}
